$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")


# --- Schedule sheet ---
$schedule.Range("A2").Value = 46079.04166666666
$schedule.Range("C2").Value = 5
$schedule.Range("D2").Value = 18.9
$schedule.Range("E2").Value = 549.61749375
$schedule.Range("F2").Value = 29.08029067460318
$schedule.Range("A4").Value = 46079.91666666666
$schedule.Range("C4").Value = 6
$schedule.Range("D4").Value = 22.68
$schedule.Range("E4").Value = 658.75133025
$schedule.Range("F4").Value = 29.04547311507937
$schedule.Range("E5").Value = 217.2849315
$schedule.Range("F5").Value = 7.185348263888888

# --- Detailed sheet ---
$detailed.Range("E3").Value = "OFF"
$detailed.Range("B38").Value = 46.56941
$detailed.Range("B39").Value = 61.07767
$detailed.Range("B40").Value = 64.14136999999999
$detailed.Range("C40").Value = "historical"
$detailed.Range("B41").Value = 63.05827
$detailed.Range("C41").Value = "historical"
$detailed.Range("B42").Value = 62.63629
$detailed.Range("C42").Value = "historical"
$detailed.Range("B43").Value = 61.75209
$detailed.Range("C43").Value = "historical"
$detailed.Range("B44").Value = 61.83105
$detailed.Range("C44").Value = "historical"
$detailed.Range("B45").Value = 65
$detailed.Range("C45").Value = "historical"
$detailed.Range("B46").Value = 57.06
$detailed.Range("C46").Value = "historical"
$detailed.Range("E46").Value = "ON"
$detailed.Range("B47").Value = 59.88233
$detailed.Range("C47").Value = "historical"
$detailed.Range("B48").Value = 57.31
$detailed.Range("C48").Value = "historical"
$detailed.Range("B49").Value = 57.06009
$detailed.Range("C49").Value = "historical"
$detailed.Range("B50").Value = 57.06
$detailed.Range("B51").Value = 56.98
$detailed.Range("B52").Value = 55.86599
$detailed.Range("B53").Value = 54.51572
$detailed.Range("B54").Value = 52.62802
$detailed.Range("B55").Value = 55.77422
$detailed.Range("B56").Value = 56.49206
$detailed.Range("B57").Value = 55.01396
$detailed.Range("B58").Value = 56.98
$detailed.Range("B59").Value = 57.06003
$detailed.Range("B60").Value = 57.53513
$detailed.Range("B61").Value = 59.43986
$detailed.Range("B62").Value = 65
$detailed.Range("B63").Value = 77.26758
$detailed.Range("B64").Value = 75.25366
$detailed.Range("B65").Value = 57.06016
$detailed.Range("B66").Value = 51.86393
$detailed.Range("B68").Value = 35.85013
$detailed.Range("B69").Value = 13.04753
$detailed.Range("B70").Value = 32.807
$detailed.Range("B71").Value = 8.657159999999999
$detailed.Range("B72").Value = -0.04679
$detailed.Range("B73").Value = 0.72986
$detailed.Range("B76").Value = 34.01368
$detailed.Range("B78").Value = -3.76
$detailed.Range("B79").Value = 8.153840000000001
$detailed.Range("B80").Value = 0.73
$detailed.Range("B81").Value = 0.73
$detailed.Range("B83").Value = 51.48038
$detailed.Range("B84").Value = 49.41504
$detailed.Range("B85").Value = 37.95
$detailed.Range("B86").Value = 50.46669
$detailed.Range("B87").Value = 64.89
$detailed.Range("B88").Value = 75.34499
$detailed.Range("B89").Value = 73.34413000000001
$detailed.Range("B90").Value = 74.22476
$detailed.Range("B91").Value = 78
$detailed.Range("B92").Value = 71.40000000000001
$detailed.Range("B93").Value = 65
$detailed.Range("B94").Value = 64.89
$detailed.Range("B95").Value = 59.38176
